$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated coin price/volume data (cryptos.xlsx symbol list refresh).
# D-column cells are numeric-looking values stored as text; E-column cells
# are plain text labels. For D-column we force NumberFormat "@" before
# assigning so Excel keeps the exact text representation (no numeric
# coercion / precision loss, trailing zeros preserved) instead of turning
# the value into a floating point number.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "248.27"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "21.75"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.364"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.408"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8154"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9530"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1411"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07599"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03194"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03051"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09303"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.560"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001610"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04707"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005770"

$ws.Range("E18").Value = "17OneONEWorstin24h"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006379"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.005072"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.001034"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.747"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.146"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3252"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0003099"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03943"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006969"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1062"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003400"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008609"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005819"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0005500"

$ws.Range("E47").Value = "46ACDXExchangeACXT"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.6799"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1660"
